# horaire: passer billeterie cours 1
#
# Move the "[Billeterie](billeterie.md)" link from the "niveaux de
# priorites / types d'utilisateurs" cell (week 1, course 2 -> C3) into the
# "Plan de cours" cell (week 1, course 1 -> C2) on the "modele" sheet.
# groupe1 / groupe2 pull C2/C3 via formulas (=modele!C2 / =modele!C3), so
# they recalculate automatically.

$wb = $excel.ActiveWorkbook
$modele = $wb.Worksheets.Item("modele")

# --- C3 (course 2): drop Billeterie, keep the other two links --------
# Leading "'" reproduces the quote-prefix + wrap-text style (xf index 7)
# that this cell picks up in the target workbook.
$modele.Range("C3").Value = "'[niveaux de priorités](niveaux-priorites.md) <br/> [Types d'utilisateurs](types-utilisateurs.md)"

# --- C2 (course 1): append the Billeterie link -----------------------
$modele.Range("C2").Value = "Plan de cours <br/>[Résolution de problèmes](resolution-probleme.md)<br/>[Introduction apprentissage par problème](apprentissage-probleme.md) <br/>[Formation des équipes pour la session](formation-equipes.md)<br/>[Billeterie](billeterie.md) "

# The longer C2 text now needs a taller row to stay fully visible.
$modele.Rows.Item(2).RowHeight = 119

# --- Active tab / selection bookkeeping -------------------------------
# Previously groupe2 was the active tab (selection F6); now modele is
# active (selection C3) and groupe2's selection moves to D3.
$groupe2 = $wb.Worksheets.Item("groupe2")
$groupe2.Activate()
$groupe2.Range("D3").Select()

$modele.Activate()
$modele.Range("C3").Select()
